# Fix a minor error on slide 5: "down key" -> "right key" plus a small
# repositioning (shift up/left) of the four numbered step captions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# The four caption text boxes live inside the big outer group ("그룹 169"),
# which is the 2nd top-level shape on the slide.
$grp = $s.Shapes.Item(2)

# Locate the four shapes we need to touch by their (stable) shape Id,
# since GroupItems indices enumerate every nested shape in the group.
$textBox25 = $null
$textBox34 = $null
$textBox61 = $null
$textBox63 = $null
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $item = $grp.GroupItems.Item($i)
    if ($item.Id -eq 26) { $textBox25 = $item }
    elseif ($item.Id -eq 35) { $textBox34 = $item }
    elseif ($item.Id -eq 62) { $textBox61 = $item }
    elseif ($item.Id -eq 64) { $textBox63 = $item }
}

# --- Reposition the four caption boxes slightly (up & to the left) ---
# NOTE: Left/Top are stored as single-precision (float32) points and then
# multiplied by 12700 to get EMUs, with truncation (not rounding) applied.
# Adding half an EMU (in points) compensates for that truncation so the
# resulting EMU values land exactly on target.
$emuEpsilon = 0.5 / 12700

$textBox25.Left = 1484404 / 12700 + $emuEpsilon
$textBox25.Top = 841251 / 12700 + $emuEpsilon

$textBox34.Left = 6305714 / 12700 + $emuEpsilon
$textBox34.Top = 846979 / 12700 + $emuEpsilon

$textBox61.Left = 1484404 / 12700 + $emuEpsilon
$textBox61.Top = 3982347 / 12700 + $emuEpsilon

$textBox63.Left = 6310020 / 12700 + $emuEpsilon
$textBox63.Top = 3907858 / 12700 + $emuEpsilon

# --- Fix the text of step 1: "down key" -> "right key" ---
$tr = $textBox25.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)
$fullText = $para1.Text
$idx = $fullText.IndexOf("down key ")
$sub = $tr.Characters($idx + 1, 9)
$sub.Text = "right key "
